# Apply "break out stock.yaml completed" edit:
#  - Convert E79 (bsecode) from text to a true number
#  - Append a new row 80 with another NMDC screener hit

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10per change")

# E79 was stored as text ("526371"); make it numeric like the other bsecode cells.
$ws.Range("E79").Value = 526371

# New row 80
$ws.Range("A80").Value = "27/06/2024 06:44:40"
$ws.Range("B80").Value = 1
$ws.Range("C80").Value = "NMDC"
$ws.Range("D80").Value = "Nmdc Limited"
$ws.Range("E80").NumberFormat = "@"
$ws.Range("E80").Value = "526371"
$ws.Range("F80").Value = -1.06
$ws.Range("G80").Value = 246.6
$ws.Range("H80").Value = 7262513
